$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- P column - "Meets Both - AND()" boolean formula ---
# Mirrors the existing L/M/N/O pattern: row 3 gets a standalone formula,
# rows 4-12 are filled down together as a shared formula group.
$ws.Range("P3").Formula = '=AND(L3,M3)'
$ws.Range("P4:P12").Formula = '=AND(L4,M4)'

# --- Row 13 - "Average (by Cell)" ---
$ws.Range("C13").Formula = '=(C3+C4+C5+C6+C7+C8+C9+C10+C11+C12)/10'
$ws.Range("D13:E13").Formula = '=(D3+D4+D5+D6+D7+D8+D9+D10+D11+D12)/10'
$ws.Range("D13").Copy()
$ws.Range("E13").PasteSpecial(-4122)

# --- Row 14 - "Average (AVERAGE())" ---
$ws.Range("C14").Formula = '=AVERAGE(C3:C12)'
$ws.Range("D14:E14").Formula = '=AVERAGE(D3:D12)'
$ws.Range("D14").Copy()
$ws.Range("E14").PasteSpecial(-4122)

# --- Row 15 - "Total Count" ---
$ws.Range("C15").Formula = '=COUNT(C3:C12)'
$ws.Range("D15:E15").Formula = '=COUNT(D3:D12)'
$ws.Range("D15").Copy()
$ws.Range("E15").PasteSpecial(-4122)

# --- Row 16 - "Meets Goals Count" (literal criteria) ---
$ws.Range("C16").Formula = '=COUNTIF($C$3:$C$12,"<=5")'
$ws.Range("D16").Formula = '=COUNTIF($D$3:$D$12, ">= "&S4)'

# --- Row 17 - "Meets Goals Count" (dynamic / hard-coded criteria) ---
$ws.Range("C17").Formula = '=COUNTIF(C3:C12,"<="&S3)'
$ws.Range("D17").Formula = '=COUNTIF($D$3:$D$12, ">= 90000")'
$ws.Range("C16").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("D16").Copy()
$ws.Range("D17").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Leave the selection where the author ended up after this work.
[void]$ws.Range("O4").Select()
